# Update numeric values in Sheet1 to reflect re-run of the RandomForest
# imputation algorithm (commit message: "Update Name of Algo").
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    "C3"  = -12.4798
    "E3"  = 15.66620000000001
    "C21" = -12.3987
    "C23" = -12.3424
    "E24" = 17.20360000000001
    "C25" = -13.3122
    "D27" = -8.849300000000005
    "D31" = -8.912400000000007
    "D39" = -7.970499999999999
    "D48" = -7.467999999999997
    "D51" = -7.710800000000002
    "D52" = -7.752800000000002
    "C53" = -10.55140000000001
    "D55" = -8.391799999999996
    "D56" = -7.844799999999999
    "C57" = -14.1462
    "D57" = -8.354899999999992
    "E57" = 16.68150000000001
    "C59" = -12.6634
    "E61" = 16.496
    "C69" = -10.7137
    "E70" = 17.09660000000001
    "D73" = -7.476699999999998
    "C79" = -10.70530000000001
    "C83" = -13.9599
    "E86" = 16.68180000000001
    "D89" = -6.017400000000005
    "D90" = -8.079300000000003
    "C93" = -11.31250000000001
    "E98" = 15.6135
    "E100" = 16.90520000000001
    "E102" = 16.52179999999998
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
